$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows above row 2 (pushes current rows 2..21 down to 10..29)
$ws.Rows("2:9").Insert()
$ws.Rows("2:9").ClearFormats()

# Fill the 8 newly inserted rows (2..9) with the new data
$newTopRows = @(
    @(0.0106901414692401, -0.00335975876078, 0.0360410511493682),
    @(-0.007177666760981, -0.0487165041267871, 0.0716239511966705),
    @(-0.0007635815418325, -0.0448985956609249, 0.0595593601465225),
    @(0.09666942805051799, 0.0059559359215199, 0.0488692186772823),
    @(0.1411098688840866, 0.2434297949075698, -0.0125227374956011),
    @(-0.0161879286170005, 0.07849618047475811, 0.0746782794594764),
    @(0.0678060427308082, -0.026419922709465, 0.0497855171561241),
    @(-0.0740674138069152, -0.4193589985370636, 0.0155770638957619)
)

$r = 2
foreach ($row in $newTopRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Append 2 new rows at the bottom (rows 30 and 31)
$newBottomRows = @(
    @(-0.0224492978304624, 0.0058032199740409, 0.0675006061792373),
    @(-0.009010262787342, -0.0429132841527462, -0.06276640295982359)
)

$r = 30
foreach ($row in $newBottomRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
